$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "66.253.27"
$ws.Cells.Item(2, 5).Value = "  -0.95%  "
$ws.Cells.Item(3, 4).Value = "3.561.80"
$ws.Cells.Item(3, 5).Value = "  +2.29%  "
$ws.Cells.Item(4, 5).Value = "  -0.10%  "
$ws.Cells.Item(5, 4).Value = "'606.61"
$ws.Cells.Item(5, 5).Value = "  +0.34%  "
$ws.Cells.Item(6, 4).Value = "'144.51"
$ws.Cells.Item(6, 5).Value = "  +0.16%  "
$ws.Cells.Item(7, 4).Value = "3.560.25"
$ws.Cells.Item(7, 5).Value = "  +2.38%  "
$ws.Cells.Item(8, 5).Value = "  +0.06%  "
$ws.Cells.Item(9, 5).Value = "  +1.13%  "
$ws.Cells.Item(10, 5).Value = "  -2.85%  "
$ws.Cells.Item(11, 4).Value = "'8.03"
$ws.Cells.Item(11, 5).Value = "  +1.52%  "
$ws.Cells.Item(12, 5).Value = "  -1.35%  "
$ws.Cells.Item(13, 4).Value = "4.164.00"
$ws.Cells.Item(14, 5).Value = "  -1.76%  "
$ws.Cells.Item(15, 4).Value = "'30.28"
$ws.Cells.Item(15, 5).Value = "  -2.04%  "
$ws.Cells.Item(16, 4).Value = "3.560.72"
$ws.Cells.Item(16, 5).Value = "  +1.88%  "
$ws.Cells.Item(17, 4).Value = "66.300.74"
$ws.Cells.Item(17, 5).Value = "  -1.04%  "
$ws.Cells.Item(18, 4).Value = "'11.66"
$ws.Cells.Item(18, 5).Value = "  +9.42%  "
$ws.Cells.Item(19, 5).Value = "  -1.15%  "
$ws.Cells.Item(20, 4).Value = "'6.21"
$ws.Cells.Item(20, 5).Value = "  -0.66%  "
$ws.Cells.Item(21, 5).Value = "  -2.07%  "
$ws.Cells.Item(22, 4).Value = "'428.97"
$ws.Cells.Item(22, 5).Value = "  +0.04%  "
$ws.Cells.Item(23, 4).Value = "'0.609"
$ws.Cells.Item(23, 5).Value = "  +1.52%  "
$ws.Cells.Item(24, 4).Value = "'78.74"
$ws.Cells.Item(24, 5).Value = "  -0.84%  "
$ws.Cells.Item(25, 4).Value = "3.700.53"
$ws.Cells.Item(25, 5).Value = "  +2.12%  "
$ws.Cells.Item(26, 5).Value = "  -0.09%  "
$ws.Cells.Item(27, 5).Value = "  +4.09%  "
$ws.Cells.Item(28, 5).Value = "  -0.03%  "
$ws.Cells.Item(29, 4).Value = "'9.25"
$ws.Cells.Item(29, 5).Value = "  -4.62%  "
$ws.Cells.Item(30, 4).Value = "'2.52"
$ws.Cells.Item(30, 5).Value = "  +0.89%  "
$ws.Cells.Item(31, 4).Value = "'0.998"
$ws.Cells.Item(31, 5).Value = "  -0.37%  "
$ws.Cells.Item(32, 5).Value = "  -3.33%  "
$ws.Cells.Item(33, 5).Value = "  -3.53%  "
$ws.Cells.Item(34, 4).Value = "3.559.08"
$ws.Cells.Item(34, 5).Value = "  +2.13%  "
$ws.Cells.Item(35, 4).Value = "'25.46"
$ws.Cells.Item(35, 5).Value = "  +0.70%  "
$ws.Cells.Item(37, 4).Value = "'1.76"
$ws.Cells.Item(37, 5).Value = "  -0.35%  "
$ws.Cells.Item(38, 2).Value = "Aptos"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(38, 4).Value = "'7.87"
$ws.Cells.Item(38, 5).Value = "  -0.29%  "
$ws.Cells.Item(39, 2).Value = "NEARProtocol"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(39, 4).Value = "'5.67"
$ws.Cells.Item(39, 5).Value = "  +0.06%  "
$ws.Cells.Item(40, 4).Value = "'0.999"
$ws.Cells.Item(40, 5).Value = "  -0.12%  "
$ws.Cells.Item(41, 4).Value = "'171.80"
$ws.Cells.Item(41, 5).Value = "  -0.97%  "
$ws.Cells.Item(42, 5).Value = "  -3.20%  "
$ws.Cells.Item(43, 5).Value = "  -0.24%  "
$ws.Cells.Item(44, 4).Value = "'0.895"
$ws.Cells.Item(44, 5).Value = "  +0.76%  "
$ws.Cells.Item(45, 5).Value = "  -3.32%  "
$ws.Cells.Item(46, 4).Value = "'45.80"
$ws.Cells.Item(47, 4).Value = "'1.22"
$ws.Cells.Item(47, 5).Value = "  +1.52%  "
$ws.Cells.Item(48, 4).Value = "'26.08"
$ws.Cells.Item(48, 5).Value = "  -5.09%  "
$ws.Cells.Item(49, 5).Value = "  +1.37%  "
$ws.Cells.Item(50, 5).Value = "  -1.81%  "
$ws.Cells.Item(51, 4).Value = "'0.951"
$ws.Cells.Item(51, 5).Value = "  -2.81%  "
